# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text updated from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Target File / Latest Handback File / Latest Handback DateTime populated
#    for both language sheets (zh-cn, de-de), plus new hyperlinks on the
#    "Latest Target File" column pointing at the source .md files
#  - A couple of columns widened so the new long file names are readable

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---- Status column updates (Overview + both language sheets) ----
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_overview.Range("E3").Value = $newStatus
$ws_overview.Range("F3").Value = $newStatus

$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("C3").Value = $newStatus

$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("C3").Value = $newStatus

# ---- zh-cn sheet: Latest Target File / Latest Handback File / DateTime ----
$zhcn_md_1 = "3f4e26c3-18fe-4bcb-9667-1f651d144e00.md"
$zhcn_md_2 = "79798927-7769-49c8-93be-b82ce2ba9fe0.md"
$zhcn_md_1_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dc9252c1a0c533feaac30c40dd604014de52823/e2e/3f4e26c3-18fe-4bcb-9667-1f651d144e00.md"
$zhcn_md_2_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dc9252c1a0c533feaac30c40dd604014de52823/e2e/79798927-7769-49c8-93be-b82ce2ba9fe0.md"

$ws_zhcn.Range("J2").Value = "3f4e26c3-18fe-4bcb-9667-1f651d144e00.e665a7055e6ebb7cf801bbcc67136da8e20b3280.zh-cn.xlf"
$ws_zhcn.Range("J3").Value = "79798927-7769-49c8-93be-b82ce2ba9fe0.73caec0e619da511e1279ef0080b0c9d621bc515.zh-cn.xlf"

$ws_zhcn.Range("K2").Value = "2016-11-09 07:29:58"
$ws_zhcn.Range("K3").Value = "2016-11-09 07:29:58"

# Recreate the hyperlinks in document order (A2, I2, A3, I3) so relationship
# ids land the same way Excel lays them out after adding the two new links.
$ws_zhcn.Hyperlinks.Delete()
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("A2"), $zhcn_md_1_url, "", "", $zhcn_md_1)
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I2"), $zhcn_md_1_url, "", "", $zhcn_md_1)
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("A3"), $zhcn_md_2_url, "", "", $zhcn_md_2)
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I3"), $zhcn_md_2_url, "", "", $zhcn_md_2)

# ---- de-de sheet: Latest Target File / Latest Handback File / DateTime ----
$dede_md_1_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dc9252c1a0c533feaac30c40dd604014de52823/e2e/3f4e26c3-18fe-4bcb-9667-1f651d144e00.md"
$dede_md_2_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dc9252c1a0c533feaac30c40dd604014de52823/e2e/79798927-7769-49c8-93be-b82ce2ba9fe0.md"

$ws_dede.Range("J2").Value = "3f4e26c3-18fe-4bcb-9667-1f651d144e00.e665a7055e6ebb7cf801bbcc67136da8e20b3280.de-de.xlf"
$ws_dede.Range("J3").Value = "79798927-7769-49c8-93be-b82ce2ba9fe0.73caec0e619da511e1279ef0080b0c9d621bc515.de-de.xlf"

$ws_dede.Range("K2").Value = "2016-11-09 07:30:17"
$ws_dede.Range("K3").Value = "2016-11-09 07:30:17"

$ws_dede.Hyperlinks.Delete()
$ws_dede.Hyperlinks.Add($ws_dede.Range("A2"), $dede_md_1_url, "", "", $zhcn_md_1)
$ws_dede.Hyperlinks.Add($ws_dede.Range("I2"), $dede_md_1_url, "", "", $zhcn_md_1)
$ws_dede.Hyperlinks.Add($ws_dede.Range("A3"), $dede_md_2_url, "", "", $zhcn_md_2)
$ws_dede.Hyperlinks.Add($ws_dede.Range("I3"), $dede_md_2_url, "", "", $zhcn_md_2)

# ---- Column width adjustments (to fit the longer values now shown) ----
$ws_overview.Columns.Item(5).ColumnWidth = 29.2
$ws_overview.Columns.Item(6).ColumnWidth = 29.2

$ws_zhcn.Columns.Item(3).ColumnWidth = 29.2
$ws_zhcn.Columns.Item(9).ColumnWidth = 39.15
$ws_zhcn.Columns.Item(10).ColumnWidth = 39.15

$ws_dede.Columns.Item(3).ColumnWidth = 29.2
$ws_dede.Columns.Item(9).ColumnWidth = 39.15
$ws_dede.Columns.Item(10).ColumnWidth = 39.15

Write-Output "Handback report generated"
